# Apply the "updated summary charts and summary reports including comments
# from Prof. Erhardt" edit to the Lexington-Fayette FAC Summary Report.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# 1. Year 1 label: 2002 -> 2012 (appears twice on Sheet1)
# ---------------------------------------------------------------------
$ws1.Range("C1").Value = "2012"
$ws1.Range("E7").Value = "2012"

# ---------------------------------------------------------------------
# 2. Updated factor table data (columns E & H) for rows 8-18.
#    Column F values are unchanged; only E (new "Year 1" average value)
#    and H (new ridership-effect absolute value) move.
# ---------------------------------------------------------------------
$ws1.Range("E8").Value  = 2271522
$ws1.Range("H8").Value  = -1000599.761

$ws1.Range("E9").Value  = 0.556275769
$ws1.Range("H9").Value  = -540469.358

$ws1.Range("E10").Value = 727375.5
$ws1.Range("H10").Value = 144328.33255

$ws1.Range("E11").Value = 20.95096521
$ws1.Range("H11").Value = -5758.482116

$ws1.Range("E12").Value = 3.9349
$ws1.Range("H12").Value = -219570.465721

$ws1.Range("E13").Value = 27769.93
$ws1.Range("H13").Value = -130988.64149

$ws1.Range("E14").Value = 7.48
$ws1.Range("H14").Value = -42923.04483100001

$ws1.Range("E15").Value = 3.8
$ws1.Range("H15").Value = -11605.279291

# Row 16 (Years Since Ride-hail Start): E/F stay blank, H value unchanged.
# Row 17 (Bike Share) and Row 18 (Electric Scooters): values unchanged.

# New Reporters row: H19 goes from blank to an explicit 0.
$ws1.Range("H19").Value = 0

# Total Modeled / Observed ridership rows.
$ws1.Range("E20").Value = 5508786.972
$ws1.Range("E21").Value = 5019458

# ---------------------------------------------------------------------
# 3. Formulas: "% Diff" columns switch from *100/ to a plain ratio
#    (the number format is changed to a percentage instead, see below).
# ---------------------------------------------------------------------
for ($r = 8; $r -le 21; $r++) {
    $ws1.Range("G$r").Formula = "=IFERROR((F$r-E$r)/E$r,0)"
}
for ($r = 8; $r -le 19; $r++) {
    $ws1.Range("I$r").Formula = "=IFERROR(H$r/`$E`$21,0)"
}

# ---------------------------------------------------------------------
# 4. Number formats: E/F/H columns -> 2-decimal number; G/I columns -> %
# ---------------------------------------------------------------------
$ws1.Range("E8:F18").NumberFormat = "#,##0.00"
$ws1.Range("H8:H18").NumberFormat = "#,##0.00"

$ws1.Range("E19:F19").NumberFormat = "#,##0.00"
$ws1.Range("H21").NumberFormat = "#,##0.00"
$ws1.Range("H19").NumberFormat = "#,##0.00"

$ws1.Range("E20:F20").NumberFormat = "#,##0.00"
$ws1.Range("H20").NumberFormat = "#,##0.00"

$ws1.Range("E21:F21").NumberFormat = "#,##0.00"

$ws1.Range("G8:G21").NumberFormat = "0.00%"
$ws1.Range("I8:I19").NumberFormat = "0.00%"
$ws1.Range("I21").NumberFormat = "0.00%"

# ---------------------------------------------------------------------
# 5. View changes on Sheet1: drop the A7 top-left scroll anchor and move
#    the active selection from K20 to H21.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("H21").Select()
